# Update the Top10 device table with new names/values (rows 2-11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("1#提升泵组", 51929),
    @("3#提升泵组", 20616),
    @("臭氧发生器", 7680),
    @("2#提升泵组", 6250),
    @("脱水机", 2295),
    @("柱塞泵", 1224),
    @("尾气破坏", 1194),
    @("大反冲洗泵", 672),
    @("空压机", 671.328),
    @("冷却塔循环泵", 540)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row++
}
